$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 18 data rows (rows 2-19). The new data set has 20
# rows (rows 2-21), so insert two more rows at the bottom, copying
# formatting (number formats etc.) down from the row above so the new rows
# keep the same styles as the existing data rows.
$ws.Rows(20).Insert(-4121, 0)
$ws.Rows(21).Insert(-4121, 0)

# New data: rates of 100 Japanese Yen (JPY), replacing the old data for
# 1 US Dollar (USD).
# Columns: nominal (A), date serial (B), rate (C), currency name (D)
$data = @(
    @(100, 44597, 66.0766),
    @(100, 44596, 66.8412),
    @(100, 44595, 66.7786),
    @(100, 44594, 67.0552),
    @(100, 44593, 67.1523),
    @(100, 44590, 67.3423),
    @(100, 44589, 68.7422),
    @(100, 44588, 69.2033),
    @(100, 44587, 69.0844),
    @(100, 44586, 68.0460),
    @(100, 44583, 67.3165),
    @(100, 44582, 66.8510),
    @(100, 44581, 67.1967),
    @(100, 44580, 66.5429),
    @(100, 44579, 66.5067),
    @(100, 44576, 66.5935),
    @(100, 44575, 65.1624),
    @(100, 44574, 64.5905),
    @(100, 44573, 64.9360),
    @(100, 44572, 64.9168)
)

$currencyName = "Японская иена"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $currencyName
}
